# doudizhu_layout.xlsx - "sound finished changed data storage"
#
# Re-shapes the data-storage layout on the "斗地主" sheet:
#  - the "name_three"/"avatar_three" anchor cells move from column N to
#    column O (their 2-wide merges shrink to 1-wide, freeing column N)
#  - the "avatar_two"/"avatar_one" anchors lose their extra merged column
#    (A2:B3 -> A2:A3, A13:B14 -> A13:A14) and the now-standalone B/A15
#    cells pick up the plain "vertical-center" style instead of
#    "center/center"
#  - the 斗地主 sheet becomes the active/selected sheet (with A13:A14
#    selected), and the login sheet is no longer the active one

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "斗地主"

# --- Unmerge the cells whose merge shape changes -----------------------
$ws.Range("A1:B1").UnMerge()
$ws.Range("A2:B3").UnMerge()
$ws.Range("N1:O1").UnMerge()
$ws.Range("N2:O3").UnMerge()
$ws.Range("A13:B14").UnMerge()
$ws.Range("A15:B15").UnMerge()

# --- Move the label values that shift from column N to column O --------
$ws.Range("O1").Value = $ws.Range("N1").Value2
$ws.Range("N1").Value = $null
$ws.Range("N1").ClearFormats()

$ws.Range("O2").Value = $ws.Range("N2").Value2
$ws.Range("N2").Value = $null
$ws.Range("N2").ClearFormats()

# --- Re-merge the cells that are still merged, but narrower ------------
$ws.Range("A2:A3").Merge()
$ws.Range("O2:O3").Merge()
$ws.Range("A13:A14").Merge()
# (A1:B1 and A15:B15 stay unmerged / split)

# --- Style updates: plain "center" -> "vertical-center only" -----------
# (xlGeneral = 1, xlCenter = -4108)
$cellsToVCenterOnly = @("A1", "B1", "O1", "B2", "N3", "B3", "B13", "B14", "A15", "B15")
foreach ($addr in $cellsToVCenterOnly) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = 1
    $c.VerticalAlignment = -4108
}

# --- Activate the 斗地主 sheet and set its selection --------------------
$ws.Select()
$ws.Range("A13:A14").Select()
